# Update the division-problem worksheet table: replace each problem's
# text with the new problem, cell by cell. We scope Find.Execute to each
# cell's own Range (and use wdReplaceOne = 1) so that duplicate problem
# text (e.g. "376÷3=" appears twice in the original) is only replaced in
# the specific cell intended, not document-wide.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellRange = $t.Cell(1, 1).Range
$cellRange.Find.Execute("204÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "360÷5=", 1) | Out-Null

$cellRange = $t.Cell(1, 2).Range
$cellRange.Find.Execute("195÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "995÷8=", 1) | Out-Null

$cellRange = $t.Cell(1, 3).Range
$cellRange.Find.Execute("238÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "936÷6=", 1) | Out-Null

$cellRange = $t.Cell(1, 4).Range
$cellRange.Find.Execute("319÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "807÷9=", 1) | Out-Null

$cellRange = $t.Cell(1, 5).Range
$cellRange.Find.Execute("264÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "864÷6=", 1) | Out-Null

$cellRange = $t.Cell(5, 1).Range
$cellRange.Find.Execute("911÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "352÷6=", 1) | Out-Null

$cellRange = $t.Cell(5, 2).Range
$cellRange.Find.Execute("867÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "949÷4=", 1) | Out-Null

$cellRange = $t.Cell(5, 3).Range
$cellRange.Find.Execute("201÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "643÷6=", 1) | Out-Null

$cellRange = $t.Cell(5, 4).Range
$cellRange.Find.Execute("376÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "158÷2=", 1) | Out-Null

$cellRange = $t.Cell(5, 5).Range
$cellRange.Find.Execute("429÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "265÷4=", 1) | Out-Null

$cellRange = $t.Cell(9, 1).Range
$cellRange.Find.Execute("463÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "511÷6=", 1) | Out-Null

$cellRange = $t.Cell(9, 2).Range
$cellRange.Find.Execute("745÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "424÷6=", 1) | Out-Null

$cellRange = $t.Cell(9, 3).Range
$cellRange.Find.Execute("344÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "538÷8=", 1) | Out-Null

$cellRange = $t.Cell(9, 4).Range
$cellRange.Find.Execute("314÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "523÷5=", 1) | Out-Null

$cellRange = $t.Cell(9, 5).Range
$cellRange.Find.Execute("620÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "791÷2=", 1) | Out-Null

$cellRange = $t.Cell(13, 1).Range
$cellRange.Find.Execute("976÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "659÷6=", 1) | Out-Null

$cellRange = $t.Cell(13, 2).Range
$cellRange.Find.Execute("581÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "366÷7=", 1) | Out-Null

$cellRange = $t.Cell(13, 3).Range
$cellRange.Find.Execute("454÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "310÷4=", 1) | Out-Null

$cellRange = $t.Cell(13, 4).Range
$cellRange.Find.Execute("382÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "840÷7=", 1) | Out-Null

$cellRange = $t.Cell(13, 5).Range
$cellRange.Find.Execute("243÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "636÷6=", 1) | Out-Null

$cellRange = $t.Cell(17, 1).Range
$cellRange.Find.Execute("939÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "726÷2=", 1) | Out-Null

$cellRange = $t.Cell(17, 2).Range
$cellRange.Find.Execute("376÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "121÷3=", 1) | Out-Null

$cellRange = $t.Cell(17, 3).Range
$cellRange.Find.Execute("558÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "113÷3=", 1) | Out-Null

$cellRange = $t.Cell(17, 4).Range
$cellRange.Find.Execute("993÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "163÷4=", 1) | Out-Null

$cellRange = $t.Cell(17, 5).Range
$cellRange.Find.Execute("632÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "449÷3=", 1) | Out-Null
